$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '62.612.55'
    'E2' = '  +2.43%  '
    'D3' = '2.961.23'
    'E3' = '  +1.12%  '
    'E4' = '  -0.08%  '
    'D5' = '595.94'
    'E5' = '  +0.66%  '
    'D6' = '145.19'
    'E6' = '  -0.08%  '
    'E7' = '  -0.05%  '
    'D8' = '2.957.00'
    'E8' = '  +0.99%  '
    'E9' = '  +0.14%  '
    'D10' = '7.35'
    'E10' = '  +5.11%  '
    'D11' = '0.143'
    'E11' = '  -0.59%  '
    'E12' = '  +1.43%  '
    'D13' = '0.0000233'
    'E13' = '  +3.14%  '
    'D14' = '33.37'
    'E14' = '  -1.15%  '
    'E15' = '  -0.39%  '
    'D16' = '3.446.97'
    'E16' = '  +0.91%  '
    'D17' = '62.479.02'
    'E17' = '  +2.33%  '
    'D18' = '6.71'
    'E18' = '  -0.12%  '
    'D19' = '2.957.00'
    'E19' = '  +0.84%  '
    'D20' = '441.35'
    'E20' = '  +1.21%  '
    'D21' = '13.43'
    'E21' = '  -0.06%  '
    'D22' = '0.672'
    'E22' = '  -1.07%  '
    'D23' = '7.10'
    'E23' = '  -0.16%  '
    'D24' = '81.82'
    'E24' = '  +0.37%  '
    'E25' = '  +0.17%  '
    'D26' = '11.99'
    'E26' = '  +0.84%  '
    'E27' = '  -3.64%  '
    'E28' = '  -0.08%  '
    'E29' = '  -0.07%  '
    'E30' = '  +0.90%  '
    'E31' = '  -5.99%  '
    'D32' = '26.60'
    'E32' = '  -0.26%  '
    'E33' = '  -2.33%  '
    'D35' = '0.0₃0876'
    'E35' = '  +0.97%  '
    'D36' = '0.991'
    'E36' = '  -2.00%  '
    'D37' = '5.63'
    'E37' = '  -0.14%  '
    'E38' = '  +3.09%  '
    'D39' = '49.60'
    'E39' = '  -0.25%  '
    'D40' = '2.92'
    'E40' = '  -2.83%  '
    'D41' = '8.56'
    'E41' = '  -0.43%  '
    'E42' = '  -4.66%  '
    'D43' = '0.282'
    'E43' = '  -1.82%  '
    'D44' = '39.10'
    'E44' = '  -7.36%  '
    'D45' = '2.717.20'
    'E45' = '  +1.02%  '
    'D46' = '135.47'
    'E46' = '  +1.53%  '
    'D47' = '0.0340'
    'E47' = '  -2.23%  '
    'D48' = '362.62'
    'E48' = '  -3.80%  '
    'E50' = '  -0.39%  '
    'D51' = '22.92'
    'E51' = '  -4.85%  '
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    if ($key[0] -eq "D") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$key]
}

Write-Output "done"
